# Edit script: update adj_close (column G) series with recomputed dividend-adjusted
# values, normalize the most-recent-row timestamp, and append a new intraday row
# reflecting the latest fetched quote (as produced by the R data-refresh script).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full recomputed adj_close series for rows 2..286 (in row order).
$adjClose = @(
    "8.46446228027344", "8.38016510009766", "8.40991687774658", "8.3305778503418", "8.35041332244873", "8.23140525817871", "8.17190074920654", "8.3305778503418", "8.3305778503418", "8.29090881347656",
    "8.25123882293701", "8.13223075866699", "8.15206623077393", "7.97355365753174", "8.05289268493652", "8.11239719390869", "8.09256172180176", "8.09256172180176", "8.13223075866699", "8.17190074920654",
    "8.17190074920654", "7.91404962539673", "7.73553705215454", "7.45785093307495", "7.93388414382935", "8.09256172180176", "8.09256172180176", "7.91404962539673", "7.91404962539673", "7.83471059799194",
    "7.83471059799194", "8.21156978607178", "8.31074333190918", "8.23140525817871", "8.31074333190918", "8.23140525817871", "8.3305778503418", "8.628098487854", "8.44958686828613", "8.42975234985352",
    "8.3305778503418", "8.39008235931396", "8.6082649230957", "8.72727298736572", "8.76694202423096", "8.628098487854", "8.9851245880127", "9.32231330871582", "9.63966941833496", "9.65950393676758",
    "9.61983394622803", "9.12396621704102", "8.94545459747314", "9.02479362487793", "9.14380168914795", "8.90578460693359", "8.82644557952881", "8.84628105163574", "8.72727298736572", "8.58842945098877",
    "8.42975234985352", "8.42975234985352", "8.48925685882568", "8.31074333190918", "8.23140525817871", "8.05289268493652", "8.13223075866699", "8.09256172180176", "8.01322269439697", "7.93388414382935",
    "8.09256172180176", "8.21156978607178", "8.21156978607178", "8.23140525817871", "7.89421463012695", "7.81487607955933", "7.85454559326172", "7.69586801528931", "7.83471059799194", "7.79504108428955",
    "7.79504108428955", "7.93388414382935", "7.93388414382935", "7.97355365753174", "7.97355365753174", "7.95371913909912", "7.97355365753174", "7.93388414382935", "7.93388414382935", "7.77520656585693",
    "8.03305816650391", "8.21156978607178", "8.03305816650391", "7.95371913909912", "7.95371913909912", "8.03305816650391", "8.37024688720703", "8.17190074920654", "8.17190074920654", "7.99338865280151",
    "8.13223075866699", "8.13223075866699", "7.93388414382935", "8.09256172180176", "8.01322269439697", "7.83471059799194", "7.83471059799194", "7.81487607955933", "7.61652851104736", "7.61652851104736",
    "7.75537204742432", "7.57685899734497", "7.57685899734497", "7.43801641464233", "7.47768592834473", "8.11239719390869", "8.17190074920654", "8.23140525817871", "8.23140525817871", "8.05289268493652",
    "8.05289268493652", "8.56859493255615", "8.19173526763916", "8.31074333190918", "8.13223075866699", "8.13223075866699", "8.13223075866699", "7.97355365753174", "7.97355365753174", "7.73553705215454",
    "7.73553705215454", "7.73553705215454", "7.89421463012695", "8.09256172180176", "8.09256172180176", "7.83471059799194", "7.83471059799194", "7.83471059799194", "7.81487607955933", "7.97355365753174",
    "7.97355365753174", "7.81487607955933", "7.63636350631714", "7.73553705215454", "7.73553705215454", "7.73553705215454", "7.87438011169434", "7.87438011169434", "7.57685899734497", "7.65619802474976",
    "7.49752044677734", "7.43801641464233", "7.59669399261475", "7.73553705215454", "7.73553705215454", "7.71570253372192", "7.43801641464233", "7.12066078186035", "7.04132223129272", "7.00165271759033",
    "7.00165271759033", "6.84297513961792", "6.94214868545532", "6.74380159378052", "7.10082626342773", "6.74380159378052", "6.54545450210571", "6.94214868545532", "7.14049577713013", "7.12066078186035",
    "7.10082626342773", "7.51735496520996", "7.23966932296753", "7.23966932296753", "7.23966932296753", "7.23966932296753", "7.21983480453491", "7.21983480453491", "6.96198320388794", "6.94214868545532",
    "6.86280965805054", "6.88264465332031", "7.14049577713013", "7.14049577713013", "6.94214868545532", "6.90247917175293", "6.94214868545532", "6.94214868545532", "6.92231416702271", "7.18016481399536",
    "7.33884286880493", "7.33884286880493", "7.33884286880493", "7.18016481399536", "7.37851238250732", "7.71570253372192", "7.53718996047974", "7.41818189620972", "7.2595043182373", "6.92231416702271",
    "6.66446256637573", "6.66446256637573", "6.66446256637573", "6.52561950683594", "6.3867769241333", "6.40661144256592", "6.42644596099854", "6.84297513961792", "6.54545450210571", "6.44628095626831",
    "6.34710741043091", "6.20826435089111", "6.08925580978394", "6.02975177764893", "6.02975177764893", "6.12892532348633", "6.12892532348633", "6.0495867729187", "6.0495867729187", "6.26776885986328",
    "6.12892532348633", "6.12892532348633", "6.12892532348633", "6.0495867729187", "5.9504132270813", "5.8512396812439", "5.83140516281128", "5.73223161697388", "5.75206613540649", "5.75206613540649",
    "5.8512396812439", "5.8512396812439", "5.9504132270813", "5.8512396812439", "5.8512396812439", "5.70247936248779", "5.60330581665039", "5.75206613540649", "5.8016529083252", "5.75206613540649",
    "5.75206613540649", "5.70247936248779", "6.24793386459351", "6.14875984191895", "6", "6", "6", "5.8512396812439", "5.70247936248779", "5.45454549789429",
    "5.30578470230103", "5.15702438354492", "5.00826454162598", "5.30578470230103", "5.30578470230103", "5.40495824813843", "5.50413227081299", "5.45454549789429", "5.20661163330078", "4.95867776870728",
    "4.95867776870728", "5.35537195205688", "5.25619840621948", "5.20661163330078", "5.30578470230103", "5.35537195205688", "5.55371904373169", "5.75206613540649", "5.60330581665039", "5.55371904373169",
    "5.55371904373169", "5.75206613540649", "5.65289211273193", "5.75206613540649", "6", "6", "5.9504132270813", "5.9504132270813", "5.9504132270813", "5.75206613540649",
    "6", "6.09999990463257", "6.25", "6.25", "6"
)

for ($i = 0; $i -lt $adjClose.Length; $i++) {
    $row = $i + 2
    # Values are written with a leading quote so Excel stores them as text
    # (matching the workbook's existing convention of keeping adj_close as
    # shared-string entries rather than numeric cells).
    $ws.Cells.Item($row, 7).Value = "'" + $adjClose[$i]
}

# Row 285 (2024-06-10) is no longer the "live" quote: normalize its timestamp
# to the standard 07:00 daily snapshot time used by all other rows.
$ws.Range("A285").Value = 45453.2916666667

# Append the new latest (intraday) row fetched by the script.
$newRow = 286
$ws.Range("A285").Copy($ws.Range("A286"))
$ws.Cells.Item($newRow, 1).Value = 45455.4667939815
$ws.Cells.Item($newRow, 2).Value = 400
$ws.Cells.Item($newRow, 3).Value = 6
$ws.Cells.Item($newRow, 4).Value = 6
$ws.Cells.Item($newRow, 5).Value = 6
$ws.Cells.Item($newRow, 6).Value = 6
$ws.Cells.Item($newRow, 7).Value = "'6"
$ws.Cells.Item($newRow, 8).Value = "ECMB.MI"
